# Generate Report for Handoff
# b.md is now ready for handoff (a new handback .xlf for it has been
# generated), so the "Status" / "Latest Handoff Datetime" data for the
# b.md row moves forward, and the hyperlink that used to (incorrectly)
# point at a.*.xlf's handoff file now reflects b.*.xlf.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md summary row.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-33-13 14:33:54"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) status/handoff file/datetime change, and the
# hyperlink display text for the handoff file needs to follow suit.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-13 14:33:50"

$zhcnLinks = @(
    @{ Ref = "A2"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/a.md"; Display = "a.md" },
    @{ Ref = "B2"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/a.md"; Display = ".md" },
    @{ Ref = "D2"; Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4a97c7b3cafbc1ec2f773e166407ecea610b27a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "F2"; Target = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/621601657a61dc2f1ddcda499b450cf24e77bfd5/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G2"; Target = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/03d41d3d2564a537036dc6f60c6d1dd2238f57ab/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A3"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/b.md"; Display = "b.md" },
    @{ Ref = "B3"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/b.md"; Display = ".md" },
    @{ Ref = "D3"; Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4a97c7b3cafbc1ec2f773e166407ecea610b27a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" },
    @{ Ref = "F3"; Target = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/621601657a61dc2f1ddcda499b450cf24e77bfd5/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G3"; Target = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/03d41d3d2564a537036dc6f60c6d1dd2238f57ab/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" }
)

$zhcn.Hyperlinks.Delete()
foreach ($link in $zhcnLinks) {
    $zhcn.Hyperlinks.Add($zhcn.Range($link.Ref), $link.Target, "", "", $link.Display)
}

# ---------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-13 14:33:54"

$dedeLinks = @(
    @{ Ref = "A2"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/a.md"; Display = "a.md" },
    @{ Ref = "B2"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/a.md"; Display = ".md" },
    @{ Ref = "D2"; Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f00ef18760a30ab84e9c4529bee5ff95a629b919/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "F2"; Target = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c99126fd5f6679dc0a2105bf395b1000a843990a/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G2"; Target = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7a506df38fda5e380040b93615a36f1d07f1ed9d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A3"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/b.md"; Display = "b.md" },
    @{ Ref = "B3"; Target = "https://github.com/OpenLocalizationTest/oltest/blob/b106d6463bed4a325df1f019d5537801f13a865b/e2e/b.md"; Display = ".md" },
    @{ Ref = "D3"; Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f00ef18760a30ab84e9c4529bee5ff95a629b919/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" },
    @{ Ref = "F3"; Target = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c99126fd5f6679dc0a2105bf395b1000a843990a/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G3"; Target = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7a506df38fda5e380040b93615a36f1d07f1ed9d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" }
)

$dede.Hyperlinks.Delete()
foreach ($link in $dedeLinks) {
    $dede.Hyperlinks.Add($dede.Range($link.Ref), $link.Target, "", "", $link.Display)
}

Write-Output "Report regenerated for handoff."
